# Update view-count ("F" column) figures on the three sheets that carry
# event data: 展览 (sheet 1), 演出 (sheet 2) and 全部类型 (sheet 4).
# 本地生活 (sheet 3) is unaffected.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# 展览 sheet (sheet1)
$ws1.Range("F3").Value  = 1800
$ws1.Range("F8").Value  = 348
$ws1.Range("F9").Value  = 1755
$ws1.Range("F11").Value = 1432
$ws1.Range("F12").Value = 819
$ws1.Range("F14").Value = 690
$ws1.Range("F15").Value = 12875
$ws1.Range("F16").Value = 12854
$ws1.Range("F18").Value = 748
$ws1.Range("F22").Value = 585
$ws1.Range("F26").Value = 16
$ws1.Range("F27").Value = 85
$ws1.Range("F28").Value = 257
$ws1.Range("F29").Value = 684

# 演出 sheet (sheet2)
$ws2.Range("F7").Value = 12

# 全部类型 sheet (sheet4)
$ws4.Range("F5").Value  = 1800
$ws4.Range("F13").Value = 348
$ws4.Range("F14").Value = 1755
$ws4.Range("F16").Value = 1432
$ws4.Range("F17").Value = 819
$ws4.Range("F20").Value = 690
$ws4.Range("F21").Value = 12875
$ws4.Range("F22").Value = 12854
$ws4.Range("F24").Value = 748
$ws4.Range("F28").Value = 585
$ws4.Range("F30").Value = 12
$ws4.Range("F35").Value = 16
$ws4.Range("F37").Value = 85
$ws4.Range("F38").Value = 257
$ws4.Range("F39").Value = 684
